$d = $word.ActiveDocument

# Locate paragraph 32, "I've worked on numerous projects ranging from robotics, games, websites and medical equipment; "
# which is the first of the 7 paragraphs (32-38) that get restructured into 10 new paragraphs
# forming the new "Let's work together..." / "Through best practices..." portfolio section.
$startPara = $d.Paragraphs.Item(32)
$endPara = $d.Paragraphs.Item(38)

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Individually and within a team, </w:t></w:r><w:r><w:t xml:space="preserve">I’ve worked on numerous projects ranging from robotics, games, websites, and medical equipment; as a </w:t></w:r><w:r><w:t>result,</w:t></w:r><w:r><w:t xml:space="preserve"> I</w:t></w:r><w:r><w:t>’ve</w:t></w:r><w:r><w:t xml:space="preserve"> learned to design and implement software based on technical requirements set up by the stakeholder.</w:t></w:r><w:r><w:t xml:space="preserve"> In addition,</w:t></w:r><w:r><w:t xml:space="preserve"> I take great strides to learn and apply new skills to develop innovative solutions to </w:t></w:r><w:r><w:t xml:space="preserve">solve </w:t></w:r><w:r><w:t>complex challenges.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">At heart, I’m a passionate problem solver who find greats thrill in meeting challenges that require thinking outside of the box. When I’m not busy trying to solve a problem, I’m </w:t></w:r><w:r><w:t xml:space="preserve">growing and </w:t></w:r><w:r><w:t xml:space="preserve">learning new </w:t></w:r><w:r><w:t>skills</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>like</w:t></w:r><w:r><w:t xml:space="preserve"> a new methodology, programming language, or even a spoken language. There is no limit to what fascinates me, there is an entire world full of new experiences and techniques to learn.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Let’s work together and coalesce our skills and experiences to solve the greatest of challenges</w:t></w:r><w:r><w:t xml:space="preserve"> whilst we learn from one another</w:t></w:r><w:r><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p/><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Through best practices and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Test Driven</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Development </w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t xml:space="preserve">I’ve worked on </w:t></w:r><w:r><w:t xml:space="preserve">numerous projects ranging from </w:t></w:r><w:r><w:t>robotics, games, websites</w:t></w:r><w:r><w:t xml:space="preserve"> and medical equipment; </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">In my free time, I look to solve problem by developing RPA scripts to automate redundant tasks in my daily life </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)
